$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the two conversion lines inside the multi-line note ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("✅ 1000 Bs = 6.99 = 27916.08 pesos", "✅ 1000 Bs = 7.03 = 28021.09 pesos")
$text = $text.Replace("✅ 27916.08 pesos = 6.97 = 968.18 Bs", "✅ 28021.09 pesos = 7.03 = 977.48 Bs")
$cell.Value2 = $text

# --- tasas sheet: update the four rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 142.25
$ws2.Range("O10").Value = 3986
$ws2.Range("N12").Value = 3987.99
$ws2.Range("O12").Value = 139.116
